# Natmi following Dr Hou advice
# Rebuild the C3-Itgax ligand-receptor table: instead of reporting only
# sending-cluster -> M2 edges with n=1 expressing cell, recompute with
# n=3 expressing cells and expand the target clusters to both ECs and M2
# for every sending cluster (ECs, FAPs, M2, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'C3'
$ws.Cells.Item(2,3).Value = 'Itgax'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 40.942832
$ws.Cells.Item(2,8).Value = 122.828496
$ws.Cells.Item(2,9).Value = 0.2583000005785167
$ws.Cells.Item(2,10).Value = 0.2583000005785167
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 15.04425566666667
$ws.Cells.Item(2,14).Value = 45.132767
$ws.Cells.Item(2,15).Value = 0.3565971854932169
$ws.Cells.Item(2,16).Value = 0.356597185493217
$ws.Cells.Item(2,17).Value = 615.9544323253814
$ws.Cells.Item(2,18).Value = 5543.589890928432
$ws.Cells.Item(2,19).Value = 0.09210905321919534
$ws.Cells.Item(2,20).Value = 0.09210905321919535

# Row 3
$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'C3'
$ws.Cells.Item(3,3).Value = 'Itgax'
$ws.Cells.Item(3,4).Value = 'M2'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 40.942832
$ws.Cells.Item(3,8).Value = 122.828496
$ws.Cells.Item(3,9).Value = 0.2583000005785167
$ws.Cells.Item(3,10).Value = 0.2583000005785167
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 27.144119
$ws.Cells.Item(3,14).Value = 81.432357
$ws.Cells.Item(3,15).Value = 0.643402814506783
$ws.Cells.Item(3,16).Value = 0.643402814506783
$ws.Cells.Item(3,17).Value = 1111.357104005008
$ws.Cells.Item(3,18).Value = 10002.21393604507
$ws.Cells.Item(3,19).Value = 0.1661909473593213
$ws.Cells.Item(3,20).Value = 0.1661909473593213

# Row 4
$ws.Cells.Item(4,1).Value = 'FAPs'
$ws.Cells.Item(4,2).Value = 'C3'
$ws.Cells.Item(4,3).Value = 'Itgax'
$ws.Cells.Item(4,4).Value = 'ECs'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 79.68771233333334
$ws.Cells.Item(4,8).Value = 239.063137
$ws.Cells.Item(4,9).Value = 0.5027335710876245
$ws.Cells.Item(4,10).Value = 0.5027335710876245
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 15.04425566666667
$ws.Cells.Item(4,14).Value = 45.132767
$ws.Cells.Item(4,15).Value = 0.3565971854932169
$ws.Cells.Item(4,16).Value = 0.356597185493217
$ws.Cells.Item(4,17).Value = 1198.842317834453
$ws.Cells.Item(4,18).Value = 10789.58086051008
$ws.Cells.Item(4,19).Value = 0.179273376502801
$ws.Cells.Item(4,20).Value = 0.179273376502801

# Row 5
$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'C3'
$ws.Cells.Item(5,3).Value = 'Itgax'
$ws.Cells.Item(5,4).Value = 'M2'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 79.68771233333334
$ws.Cells.Item(5,8).Value = 239.063137
$ws.Cells.Item(5,9).Value = 0.5027335710876245
$ws.Cells.Item(5,10).Value = 0.5027335710876245
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 27.144119
$ws.Cells.Item(5,14).Value = 81.432357
$ws.Cells.Item(5,15).Value = 0.643402814506783
$ws.Cells.Item(5,16).Value = 0.643402814506783
$ws.Cells.Item(5,17).Value = 2163.052746413768
$ws.Cells.Item(5,18).Value = 19467.47471772391
$ws.Cells.Item(5,19).Value = 0.3234601945848235
$ws.Cells.Item(5,20).Value = 0.3234601945848235

# Row 6
$ws.Cells.Item(6,1).Value = 'M2'
$ws.Cells.Item(6,2).Value = 'C3'
$ws.Cells.Item(6,3).Value = 'Itgax'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 37.53186833333334
$ws.Cells.Item(6,8).Value = 112.595605
$ws.Cells.Item(6,9).Value = 0.2367809244903433
$ws.Cells.Item(6,10).Value = 0.2367809244903433
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 15.04425566666667
$ws.Cells.Item(6,14).Value = 45.132767
$ws.Cells.Item(6,15).Value = 0.3565971854932169
$ws.Cells.Item(6,16).Value = 0.356597185493217
$ws.Cells.Item(6,17).Value = 564.6390228543372
$ws.Cells.Item(6,18).Value = 5081.751205689035
$ws.Cells.Item(6,19).Value = 0.08443541125173833
$ws.Cells.Item(6,20).Value = 0.08443541125173834

# Row 7
$ws.Cells.Item(7,1).Value = 'M2'
$ws.Cells.Item(7,2).Value = 'C3'
$ws.Cells.Item(7,3).Value = 'Itgax'
$ws.Cells.Item(7,4).Value = 'M2'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 37.53186833333334
$ws.Cells.Item(7,8).Value = 112.595605
$ws.Cells.Item(7,9).Value = 0.2367809244903433
$ws.Cells.Item(7,10).Value = 0.2367809244903433
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 27.144119
$ws.Cells.Item(7,14).Value = 81.432357
$ws.Cells.Item(7,15).Value = 0.643402814506783
$ws.Cells.Item(7,16).Value = 0.643402814506783
$ws.Cells.Item(7,17).Value = 1018.769500332332
$ws.Cells.Item(7,18).Value = 9168.925502990985
$ws.Cells.Item(7,19).Value = 0.152345513238605
$ws.Cells.Item(7,20).Value = 0.1523455132386049

# Row 8
$ws.Cells.Item(8,1).Value = 'sCs'
$ws.Cells.Item(8,2).Value = 'C3'
$ws.Cells.Item(8,3).Value = 'Itgax'
$ws.Cells.Item(8,4).Value = 'ECs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.3464216666666666
$ws.Cells.Item(8,8).Value = 1.039265
$ws.Cells.Item(8,9).Value = 0.002185503843515531
$ws.Cells.Item(8,10).Value = 0.002185503843515531
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 15.04425566666667
$ws.Cells.Item(8,14).Value = 45.132767
$ws.Cells.Item(8,15).Value = 0.3565971854932169
$ws.Cells.Item(8,16).Value = 0.356597185493217
$ws.Cells.Item(8,17).Value = 5.211656121806111
$ws.Cells.Item(8,18).Value = 46.90490509625499
$ws.Cells.Item(8,19).Value = 0.0007793445194822465
$ws.Cells.Item(8,20).Value = 0.0007793445194822466

# Row 9
$ws.Cells.Item(9,1).Value = 'sCs'
$ws.Cells.Item(9,2).Value = 'C3'
$ws.Cells.Item(9,3).Value = 'Itgax'
$ws.Cells.Item(9,4).Value = 'M2'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.3464216666666666
$ws.Cells.Item(9,8).Value = 1.039265
$ws.Cells.Item(9,9).Value = 0.002185503843515531
$ws.Cells.Item(9,10).Value = 0.002185503843515531
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 27.144119
$ws.Cells.Item(9,14).Value = 81.432357
$ws.Cells.Item(9,15).Value = 0.643402814506783
$ws.Cells.Item(9,16).Value = 0.643402814506783
$ws.Cells.Item(9,17).Value = 9.403310944178331
$ws.Cells.Item(9,18).Value = 84.62979849760498
$ws.Cells.Item(9,19).Value = 0.001406159324033285
$ws.Cells.Item(9,20).Value = 0.001406159324033285
